# Apply updated dSF (column F) values as part of a data repull / mean
# recalculation pass. Only the cells whose underlying source values changed
# are touched; everything else in the sheet is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F4"  = -4
    "F5"  = -5
    "F6"  = -4
    "F8"  = -4
    "F14" = -8
    "F15" = -7
    "F16" = -3
    "F18" = 0
    "F22" = -6
    "F23" = -1
    "F24" = 0
    "F29" = 5
    "F32" = -5
    "F34" = -4
    "F36" = 4
    "F38" = -1
    "F42" = 4
    "F43" = -4
    "F47" = 3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
